# Update cryptocurrency list: refresh Price/Volume(1h) figures and
# swap the VeChain/Stacks rows (rows 40-41) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to plain text so that
# values like "340.55" or "0.580" are not reinterpreted as numbers,
# matching how the source data is stored (inline/shared strings).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '52.221.85'
$ws.Range("E2").Value = '  +5.67%  '

$ws.Range("D3").Value = '2.790.26'
$ws.Range("E3").Value = '  +6.13%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '116.44'
$ws.Range("E5").Value = '  +4.15%  '

$ws.Range("D6").Value = '340.55'

$ws.Range("E7").Value = '  +5.66%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  +5.76%  '

$ws.Range("D10").Value = '42.04'
$ws.Range("E10").Value = '  +6.64%  '

$ws.Range("D11").Value = '0.0864'
$ws.Range("E11").Value = '  +6.56%  '

$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("E13").Value = '  +2.45%  '

$ws.Range("E14").Value = '  +0.64%  '

$ws.Range("D15").Value = '3.231.99'
$ws.Range("E15").Value = '  +6.37%  '

$ws.Range("D16").Value = '2.798.33'
$ws.Range("E16").Value = '  +6.42%  '

$ws.Range("D17").Value = '0.884'
$ws.Range("E17").Value = '  +3.64%  '

$ws.Range("D18").Value = '52.050.72'
$ws.Range("E18").Value = '  +5.39%  '

$ws.Range("E19").Value = '  +10.58%  '

$ws.Range("D20").Value = '13.32'
$ws.Range("E20").Value = '  -0.34%  '

$ws.Range("E21").Value = '  +4.84%  '

$ws.Range("E22").Value = '  +3.69%  '

$ws.Range("D23").Value = '278.04'
$ws.Range("E23").Value = '  +3.53%  '

$ws.Range("D24").Value = '70.15'
$ws.Range("E24").Value = '  +1.69%  '

$ws.Range("D25").Value = '2.76'
$ws.Range("E25").Value = '  +8.48%  '

$ws.Range("D26").Value = '26.80'
$ws.Range("E26").Value = '  +3.22%  '

$ws.Range("D28").Value = '10.25'
$ws.Range("E28").Value = '  +1.11%  '

$ws.Range("E29").Value = '  +1.19%  '

$ws.Range("E30").Value = '  +4.19%  '

$ws.Range("D31").Value = '34.93'
$ws.Range("E31").Value = '  +1.64%  '

$ws.Range("D32").Value = '50.38'
$ws.Range("E32").Value = '  +1.63%  '

$ws.Range("D33").Value = '5.74'
$ws.Range("E33").Value = '  +5.13%  '

$ws.Range("D34").Value = '0.0827'
$ws.Range("E34").Value = '  +2.64%  '

$ws.Range("E35").Value = '  +4.76%  '

$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("D37").Value = '18.95'
$ws.Range("E37").Value = '  -0.26%  '

$ws.Range("E38").Value = '  -0.19%  '

$ws.Range("E39").Value = '  +4.87%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '2.74'
$ws.Range("E40").Value = '  +28.90%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0375'
$ws.Range("E41").Value = '  +12.19%  '

$ws.Range("E42").Value = '  +4.18%  '

$ws.Range("D43").Value = '23.37'
$ws.Range("E43").Value = '  +3.03%  '

$ws.Range("D44").Value = '2.33'
$ws.Range("E44").Value = '  -1.31%  '

$ws.Range("D45").Value = '125.12'
$ws.Range("E45").Value = '  -3.72%  '

$ws.Range("D46").Value = '2.094.81'
$ws.Range("E46").Value = '  +1.62%  '

$ws.Range("D47").Value = '3.32'
$ws.Range("E47").Value = '  +1.67%  '

$ws.Range("E48").Value = '  +3.52%  '

$ws.Range("E49").Value = '  +7.34%  '

$ws.Range("D50").Value = '9.00'
$ws.Range("E50").Value = '  +1.85%  '

$ws.Range("D51").Value = '0.890'
$ws.Range("E51").Value = '  +20.26%  '
